$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOPICS")
$ws.Activate()

$ws.Range("A50").Value = "CDN"
$ws.Range("A48").Value = "Deploying Frontends on AWS -->  Cloudfront --> S3 "

$ws.Range("A46").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A50").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(48).RowHeight = 28.8

$ws.Range("Q46").Select() | Out-Null
